$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (C) column date value for rows 2-10 from 45170 (2023-09-01)
# to 45174 (2023-09-05), keeping the existing date formatting/style intact.
foreach ($row in 2..10) {
    $ws.Range("C$row").Value = 45174
}
